$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellData = @(
    @(1, 1, "88 x 23", "  2    3", "8|    |", "8|    |"),
    @(1, 2, "54 x 64", "  6    4", "5|    |", "4|    |"),
    @(1, 3, "35 x 43", "  4    3", "3|    |", "5|    |"),
    @(2, 1, "86 x 46", "  4    6", "8|    |", "6|    |"),
    @(2, 2, "79 x 42", "  4    2", "7|    |", "9|    |"),
    @(2, 3, "85 x 31", "  3    1", "8|    |", "5|    |"),
    @(3, 1, "73 x 53", "  5    3", "7|    |", "3|    |"),
    @(3, 2, "44 x 96", "  9    6", "4|    |", "4|    |"),
    @(3, 3, "15 x 43", "  4    3", "1|    |", "5|    |"),
    @(4, 1, "49 x 11", "  1    1", "4|    |", "9|    |"),
    @(4, 2, "94 x 86", "  8    6", "9|    |", "4|    |"),
    @(4, 3, "99 x 16", "  1    6", "9|    |", "9|    |"),
    @(5, 1, "45 x 72", "  7    2", "4|    |", "5|    |"),
    @(5, 2, "60 x 79", "  7    9", "6|    |", "0|    |"),
    @(5, 3, "54 x 97", "  9    7", "5|    |", "4|    |")
)

foreach ($entry in $cellData) {
    $row = $entry[0]
    $col = $entry[1]
    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 2
    $newText = $entry[2] + [char]11 + $entry[3] + [char]11 + "  ----" + [char]11 + $entry[4] + [char]11 + $entry[5]
    $r.Text = $newText
}
